$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'31.056.10"
$ws.Range("E2").Value = "'  +1.27%  "
$ws.Range("D3").Value = "'1.956.05"
$ws.Range("E3").Value = "'  -0.32%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("D5").Value = "'245.76"
$ws.Range("E5").Value = "'  -1.17%  "
$ws.Range("E6").Value = "'  +0.08%  "
$ws.Range("D7").Value = "'0.4884"
$ws.Range("E7").Value = "'  +1.23%  "
$ws.Range("E8").Value = "'  -0.07%  "
$ws.Range("D9").Value = "'0.06837"
$ws.Range("E9").Value = "'  +0.63%  "
$ws.Range("D10").Value = "'19.21"
$ws.Range("E10").Value = "'  -0.88%  "
$ws.Range("D11").Value = "'107.27"
$ws.Range("E11").Value = "'  -3.22%  "
$ws.Range("B12").Value = "'WrappedEther"
$ws.Range("C12").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.960.49"
$ws.Range("E12").Value = "'  -0.41%  "
$ws.Range("B13").Value = "'TRON"
$ws.Range("C13").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.07807"
$ws.Range("E13").Value = "'  +0.94%  "
$ws.Range("D14").Value = "'5.465"
$ws.Range("E14").Value = "'  -0.51%  "
$ws.Range("D15").Value = "'0.7017"
$ws.Range("E15").Value = "'  +1.57%  "
$ws.Range("D16").Value = "'283.51"
$ws.Range("E16").Value = "'  -3.43%  "
$ws.Range("D17").Value = "'31.063.76"
$ws.Range("E17").Value = "'  +1.24%  "
$ws.Range("D18").Value = "'13.21"
$ws.Range("E18").Value = "'  -0.56%  "
$ws.Range("D19").Value = "'0.000007708"
$ws.Range("E19").Value = "'  +0.20%  "
$ws.Range("D20").Value = "'2.196.07"
$ws.Range("E20").Value = "'  -0.91%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "'  +0.05%  "
$ws.Range("D22").Value = "'5.488"
$ws.Range("E22").Value = "'  -3.36%  "
$ws.Range("E23").Value = "'  -0.09%  "
$ws.Range("D24").Value = "'6.484"
$ws.Range("E24").Value = "'  -1.83%  "
$ws.Range("D25").Value = "'9.837"
$ws.Range("E25").Value = "'  -0.88%  "
$ws.Range("D26").Value = "'169.83"
$ws.Range("E26").Value = "'  -0.67%  "
$ws.Range("D27").Value = "'20.00"
$ws.Range("E27").Value = "'  -0.66%  "
$ws.Range("D28").Value = "'2.200"
$ws.Range("E28").Value = "'  +0.07%  "
$ws.Range("D29").Value = "'0.1058"
$ws.Range("E29").Value = "'  -1.31%  "
$ws.Range("D30").Value = "'1.410"
$ws.Range("E30").Value = "'  -2.10%  "
$ws.Range("D31").Value = "'1.586"
$ws.Range("E31").Value = "'  -1.16%  "
$ws.Range("D32").Value = "'4.610"
$ws.Range("E32").Value = "'  -1.82%  "
$ws.Range("D33").Value = "'4.452"
$ws.Range("E33").Value = "'  -0.30%  "
$ws.Range("D34").Value = "'0.04946"
$ws.Range("E34").Value = "'  -3.58%  "
$ws.Range("E35").Value = "'  -1.81%  "
$ws.Range("D36").Value = "'1.174"
$ws.Range("E36").Value = "'  -0.55%  "
$ws.Range("E37").Value = "'  -0.18%  "
$ws.Range("D38").Value = "'0.02011"
$ws.Range("D39").Value = "'2.705"
$ws.Range("E39").Value = "'  -0.25%  "
$ws.Range("D40").Value = "'6.530"
$ws.Range("E40").Value = "'  +6.43%  "
$ws.Range("D41").Value = "'2.121"
$ws.Range("E41").Value = "'  +2.64%  "
$ws.Range("D42").Value = "'74.79"
$ws.Range("E42").Value = "'  +6.49%  "
$ws.Range("D43").Value = "'0.8868"
$ws.Range("E43").Value = "'  +1.35%  "
$ws.Range("D44").Value = "'0.4470"
$ws.Range("E44").Value = "'  -0.04%  "
$ws.Range("D45").Value = "'109.46"
$ws.Range("E45").Value = "'  -1.47%  "
$ws.Range("D46").Value = "'8.174"
$ws.Range("E46").Value = "'  +10.51%  "
$ws.Range("E47").Value = "'  +0.04%  "
$ws.Range("D48").Value = "'1.001.58"
$ws.Range("E48").Value = "'  +10.98%  "
$ws.Range("E49").Value = "'  -1.22%  "
$ws.Range("D50").Value = "'9.332"
$ws.Range("E50").Value = "'  -0.39%  "
$ws.Range("D51").Value = "'0.2598"
$ws.Range("E51").Value = "'  +3.48%  "
